$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows above the current row 11 ("Job"/"Sjoerd" data),
# shifting the existing two data rows down to rows 13 and 14.
$ws.Rows.Item(11).Resize(2).Insert()

# New row 11: Frank
$ws.Cells.Item(11, 10).Value = 1
$ws.Cells.Item(11, 11).Value = "Frank"
$ws.Cells.Item(11, 12).Value = 30
$ws.Cells.Item(11, 13).Value = 8

# New row 12: Peter
$ws.Cells.Item(12, 10).Value = 2
$ws.Cells.Item(12, 11).Value = "Peter"
$ws.Cells.Item(12, 12).Value = 45
$ws.Cells.Item(12, 13).Value = 22

# Row 13 (was row 11: Job) and row 14 (was row 12: Sjoerd) keep their original
# values automatically after the insert, no changes needed there.

# New row 15: Johan, appended after the last existing data row.
$ws.Cells.Item(15, 10).Value = 5
$ws.Cells.Item(15, 11).Value = "Johan"
$ws.Cells.Item(15, 12).Value = 42
$ws.Cells.Item(15, 13).Value = 18
